$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Fix two pre-existing rounding typos that were corrected in this revision
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = 2.96
$ws.Range("E7").Value = 2.35

# ---------------------------------------------------------------------------
# 2. Insert a new data row after row 7 (Przybylski et al., 2014, Study 5)
#    This pushes the old row 8 (and everything below) down by one row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

$ws.Range("A6").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").Value = "Przybylski et al., 2014, Study 5"
$ws.Range("B8").Value = 0.03
$ws.Range("C8").Value = "[-.16, .22]"
$ws.Range("D8").Value = 109
$ws.Range("E8").Value = 2.96
$ws.Range("F8").Value = 38.12

# ---------------------------------------------------------------------------
# 3. Insert a new data row after row 18 (Tear & Nielsen, 2014, hurting
#    behavior in Tangram task). At this point in the sheet (after the first
#    insertion above) that footnote row is row 18 and the "Aggressive
#    Cognition" header is row 19, so the new row is inserted at position 19.
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Insert()

$ws.Range("A6").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null

$ws.Range("A19").Value = "Tear & Nielsen, 2014, hurting behavior in Tangram task"
$ws.Range("B19").Value = 0.01
$ws.Range("C19").Value = "[-.17, .19]"
$ws.Range("D19").Value = 120
$ws.Range("E19").Value = 3.6
$ws.Range("F19").Value = 9.01

# ---------------------------------------------------------------------------
# 4. Add a new footnote two rows below the existing footnotes (row 28),
#    leaving row 27 blank as in the reference layout.
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Tear & Nielsen, 2014; Violent & Ultraviolent conditions combined"

# ---------------------------------------------------------------------------
# 5. Cosmetic sheet-level tweaks: widen column A, update the saved selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 57.4

$ws.Range("H20").Select()
